$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 180, shifting existing rows 180..281 down to 181..282
$ws.Rows("180:180").Insert()

# Populate the newly inserted row 180 with the new data record
$ws.Range("A180").Value = 10
$ws.Range("B180").Value = "Vega Modelo de Temuco"
$ws.Range("C180").Value = "La Araucanía"
$ws.Range("D180").Value = 44806
$ws.Range("D180").NumberFormat = $ws.Range("D181").NumberFormat
$ws.Range("E180").Value = 9
$ws.Range("F180").Value = 100112039
$ws.Range("G180").Value = "Ciboulette"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 65
$ws.Range("K180").Value = 7000
$ws.Range("L180").Value = 7000
$ws.Range("M180").Value = 7000
$ws.Range("N180").Value = "$/docena de atados"
$ws.Range("O180").Value = "Provincia de Cautín"
$ws.Range("P180").Value = 2333
$ws.Range("Q180").Value = 3
$ws.Range("R180").Value = "Hortaliza"
